# Adding the changes we made on may 9th
# Prepend 9 newly captured "falling" sensor samples ahead of the existing
# series (their ax..gz readings replace C2:H10 in place), and append a new
# batch of 10 more samples (rows 22:31) continuing the timestamp sequence.
# timestamp/label columns (A/B) for the already-existing rows are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: [timestamp, label, ax, ay, az, gx, gy, gz]
$data = @{}
$data[2] = @(0, "falling", 0.2728629112243653, -0.005532175302505469, 0.17908151820302, 0.0500909499824047, 0.0342084541916847, 0.0232128798961639)
$data[3] = @(100, "falling", -0.1289855480194095, -0.0322514295578003, 0.2833916515111924, 0.0068722339347004, 0.0074830991216003, 0.0390953756868839)
$data[4] = @(200, "falling", -0.2541066646575928, -0.01140816211700431, 0.3126487381756305, -0.0138971842825412, 0.0290160998702049, 0.1440114825963974)
$data[5] = @(300, "falling", -0.1588943481445311, -0.05039391517639156, 0.1764491081237791, 0.0100792767480015, -0.022754730656743, 0.0288633834570646)
$data[6] = @(400, "falling", 0.003789019584655856, -0.07255983948707578, 0.1795819453895092, -0.0282525178045034, -0.020616702735424, -0.0172569435089826)
$data[7] = @(500, "falling", 0.0466578006744384, -0.03053182363510123, 0.2458087503910065, -0.0167987942695617, -0.0216857157647609, 0)
$data[8] = @(600, "falling", -0.2323491334915166, -0.003547763824462882, 0.2952604919672013, -0.0178678091615438, 0.0054977871477603, 0.0299323964864015)
$data[9] = @(700, "falling", -0.07431058883666952, 0.007624650001525833, 0.2976654559373855, -0.0229074470698833, 0.00534507073462, 0.0030543261673301)
$data[10] = @(800, "falling", 0.08380470275878904, 0.005895948410034098, 0.3545163981616498, -0.0394008085131645, 0.0178678091615438, 0.011148290708661)
$data[11] = @(900, "falling", -0.03897037506103547, -0.04966262578964246, 0.6067050054669385, -0.0704022198915481, 0.1944078654050827, -0.0245873257517814)
$data[12] = @(1000, "falling", -0.3774656057357791, -0.08384630084037778, 1.005563378334046, 0.3005456924438476, 0.8894197940826416, 0.086895577609539)
$data[13] = @(1100, "falling", -0.389985084533691, -0.08737320899963376, 1.122915458679199, 0.266642689704895, 1.489289522171021, 1.197143197059631)
$data[14] = @(1200, "falling", -3.129088830947885, -0.5355502605438246, 2.341251343488698, -0.5832235813140869, 0.5451972484588623, 0.808174729347229)
$data[15] = @(1300, "falling", -4.376520133018492, -0.8961358070373535, 3.528440594673157, -0.5074763298034668, -0.9810495972633362, 0.5829181671142578)
$data[16] = @(1400, "falling", -2.786781752109524, -1.147766584157945, 5.211621630191809, -0.9239336848258972, -0.2371684312820434, 0.1962404549121856)
$data[17] = @(1500, "falling", 1.074200153350837, 1.294980049133306, 4.674310564994808, -0.8011497855186462, -1.048092007637024, -0.0548251569271087)
$data[18] = @(1600, "falling", 33.04098894596114, -6.460391509532979, -13.8664929449559, 0.3634648323059082, 1.29381263256073, 1.935221076011657)
$data[19] = @(1700, "falling", 33.17409253120405, -12.13667659759519, -16.41260833740226, -4.552626132965088, -3.393509149551392, 1.527773976325989)
$data[20] = @(1800, "falling", -4.199231290817252, -1.979123908281309, 2.645077538490293, 0.7554876208305359, 1.477530360221863, -0.5998696684837341)
$data[21] = @(1900, "falling", -1.133084297180169, 1.222284126281733, 2.30235185623169, 0.3274237811565399, -0.0329867228865623, -0.4276056587696075)
$data[22] = @(2000, "falling", 1.011236310005188, -0.08591727167367935, 1.941863000392914, 0.0789543315768241, -0.0847575515508651, -0.1519527286291122)
$data[23] = @(2100, "falling", 1.760050582885739, -0.07424210608005596, 0.6869683876633619, -0.0221438650041818, 0.1513418704271316, -0.171500414609909)
$data[24] = @(2200, "falling", -0.4296665787696987, -0.05021018907427618, 0.4772178567945989, 0.0762054398655891, -0.0978911519050598, 0.366213709115982)
$data[25] = @(2300, "falling", -1.509320116043069, -0.4450684934854586, 1.374938857555391, -0.2686280012130737, 1.315498352050781, -0.0265726372599601)
$data[26] = @(2400, "falling", 2.49392051696776, -1.748301430046552, 0.9690718531608501, -0.1499674171209335, -0.5490151643753052, 0.4234823286533355)
$data[27] = @(2500, "falling", -0.5127081871032715, -0.141617327928543, -0.73384278640151, -0.09239336848258969, -0.09071348607540131, 0.0775798857212066)
$data[28] = @(2600, "falling", -0.1017783880233756, 0.3632039599120608, 0.1947979252785462, 0.0155770638957619, -0.1032362282276153, 0.09071348607540131)
$data[29] = @(2700, "falling", 0.2760831832885773, 0.05902776718139471, 0.4707315444946266, 0.4306600093841553, 0.7470881938934326, -0.1643227487802505)
$data[30] = @(2800, "falling", 0.7326052427291816, -0.04417074620723105, 0.2964785575866773, 1.219439744949341, 1.55587375164032, 0.009468411095440299)
$data[31] = @(2900, "falling", -0.5065834045410164, 1.347906202077864, 1.883451831340787, -0.3419318199157715, 2.081828832626343, -0.5186246037483215)

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}